# Update leve-profit figures (H..N columns) on the Typhon_Profits sheets.
# Generated from the recorded market-price refresh diff (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 3 - One for the Books
$ws.Range("H3").Value = 28000
$ws.Range("J3").Value = 28000
$ws.Range("L3").Value = 28000
$ws.Range("N3").Value = -28228

# row 69 - Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 1558.9286
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 1531.4814
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 4594.4442
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -6342.4442

# row 72 - Surgical Substitution (L)
$ws.Range("H72").Value = 1558.9286
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 1531.4814
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 13783.3326
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -22519.3326

# row 76 - Warding Off Temptation
$ws.Range("H76").Value = 3412.1052
$ws.Range("I76").Value = 3395.625
$ws.Range("K76").Value = 3395.625
$ws.Range("M76").Value = -3080.625

# row 79 - The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3412.1052
$ws.Range("I79").Value = 3395.625
$ws.Range("K79").Value = 3395.625
$ws.Range("M79").Value = -2303.625

# row 80 - Cleansing the Wicked Humours
$ws.Range("H80").Value = 12183972
$ws.Range("I80").Value = 380
$ws.Range("J80").Value = 16245169
$ws.Range("K80").Value = 1140
$ws.Range("L80").Value = 48735507
$ws.Range("M80").Value = -142
$ws.Range("N80").Value = -48737503

# row 83 - Washing Away the Sins (L)
$ws.Range("H83").Value = 12183972
$ws.Range("I83").Value = 380
$ws.Range("J83").Value = 16245169
$ws.Range("K83").Value = 3420
$ws.Range("L83").Value = 146206521
$ws.Range("M83").Value = 1572
$ws.Range("N83").Value = -146216505

# row 102 - Spell-rebound
$ws.Range("H102").Value = 28000
$ws.Range("J102").Value = 28000
$ws.Range("L102").Value = 28000
$ws.Range("N102").Value = -34490

# row 107 - Another Man's Ink
$ws.Range("H107").Value = 909.85
$ws.Range("I107").Value = 591.7222
$ws.Range("J107").Value = 3773
$ws.Range("K107").Value = 591.7222
$ws.Range("L107").Value = 3773
$ws.Range("M107").Value = 1328.2778
$ws.Range("N107").Value = -7613

# row 129 - Practical Command
$ws.Range("H129").Value = 250965.12
$ws.Range("J129").Value = 323762.16
$ws.Range("L129").Value = 971286.48
$ws.Range("N129").Value = -981286.48

# row 132 - Fast-forwarding Flora
$ws.Range("H132").Value = 2720.6099
$ws.Range("I132").Value = 2940.389
$ws.Range("J132").Value = 1138.2
$ws.Range("K132").Value = 8821.167000000001
$ws.Range("L132").Value = 3414.6
$ws.Range("M132").Value = -6291.167000000001
$ws.Range("N132").Value = -8474.6

# row 133 - Big Brush, Big Dreams
$ws.Range("H133").Value = 48500
$ws.Range("J133").Value = 48500
$ws.Range("L133").Value = 48500
$ws.Range("N133").Value = -58620

# row 137 - Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1224.5312
$ws.Range("I137").Value = 1283.7894
$ws.Range("J137").Value = 1137.9231
$ws.Range("K137").Value = 3851.3682
$ws.Range("L137").Value = 3413.7693
$ws.Range("M137").Value = -1301.3682
$ws.Range("N137").Value = -8513.7693

$ws = $wb.Worksheets.Item("ARM")
# row 132 - Don't Bore Me, Ore Me
$ws.Range("H132").Value = 13134
$ws.Range("I132").Value = 1310.4
$ws.Range("J132").Value = 64862.25
$ws.Range("K132").Value = 3931.2
$ws.Range("L132").Value = 194586.75
$ws.Range("M132").Value = -1401.2
$ws.Range("N132").Value = -199646.75

$ws = $wb.Worksheets.Item("BSM")
# row 86 - Through Thick and Thin
$ws.Range("H86").Value = 1710.5
$ws.Range("I86").Value = 1600.2273
$ws.Range("J86").Value = 1820.7727
$ws.Range("K86").Value = 1600.2273
$ws.Range("L86").Value = 1820.7727
$ws.Range("M86").Value = -477.2273
$ws.Range("N86").Value = -4066.7727

# row 89 - Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1710.5
$ws.Range("I89").Value = 1600.2273
$ws.Range("J89").Value = 1820.7727
$ws.Range("K89").Value = 8001.136500000001
$ws.Range("L89").Value = 9103.863499999999
$ws.Range("M89").Value = -2385.136500000001
$ws.Range("N89").Value = -20335.8635

$ws = $wb.Worksheets.Item("CRP")
# row 31 - Wall Not Found
$ws.Range("H31").Value = 4443.577
$ws.Range("I31").Value = 5154
$ws.Range("J31").Value = 4230.45
$ws.Range("K31").Value = 5154
$ws.Range("L31").Value = 4230.45
$ws.Range("M31").Value = -4859
$ws.Range("N31").Value = -4820.45

# row 34 - Armoires of the Rich and Famous
$ws.Range("H34").Value = 4443.577
$ws.Range("I34").Value = 5154
$ws.Range("J34").Value = 4230.45
$ws.Range("K34").Value = 5154
$ws.Range("L34").Value = 4230.45
$ws.Range("M34").Value = -4952
$ws.Range("N34").Value = -4634.45

# row 132 - Hull Lotta Damage
$ws.Range("H132").Value = 3393.682
$ws.Range("J132").Value = 5475.143
$ws.Range("L132").Value = 16425.429
$ws.Range("N132").Value = -21485.429

# row 134 - Wood You Be Quiet
$ws.Range("H134").Value = 1688.9375
$ws.Range("I134").Value = 1617.1538
$ws.Range("K134").Value = 4851.4614
$ws.Range("M134").Value = -2316.4614

$ws = $wb.Worksheets.Item("CUL")
# row 131 - The Mountain Steeped
$ws.Range("H131").Value = 766.4433
$ws.Range("I131").Value = 665
$ws.Range("J131").Value = 768.5789
$ws.Range("K131").Value = 1995
$ws.Range("L131").Value = 2305.7367
$ws.Range("M131").Value = 3045
$ws.Range("N131").Value = -12385.7367

$ws = $wb.Worksheets.Item("GSM")
# row 5 - Hora at Me
$ws.Range("H5").Value = 8399.933999999999
$ws.Range("I5").Value = 4333
$ws.Range("J5").Value = 9416.666999999999
$ws.Range("K5").Value = 4333
$ws.Range("L5").Value = 9416.666999999999
$ws.Range("M5").Value = -4221
$ws.Range("N5").Value = -9640.666999999999

# row 80 - Needs More Prayerbell
$ws.Range("H80").Value = 4442.857
$ws.Range("J80").Value = 4922.222
$ws.Range("L80").Value = 4922.222
$ws.Range("N80").Value = -6918.222

# row 83 - With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4442.857
$ws.Range("J83").Value = 4922.222
$ws.Range("L83").Value = 24611.11
$ws.Range("N83").Value = -34595.11

# row 132 - On Board for Lar
$ws.Range("H132").Value = 24071.262
$ws.Range("I132").Value = 1873.7059
$ws.Range("K132").Value = 5621.1177
$ws.Range("M132").Value = -3091.1177

$ws = $wb.Worksheets.Item("LTW")
# row 68 - You Could Say It's a Moving Target
$ws.Range("H68").Value = 2699.6667
$ws.Range("I68").Value = 2700
$ws.Range("J68").Value = 2699.5
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 2699.5
$ws.Range("M68").Value = -1951
$ws.Range("N68").Value = -4197.5

# row 71 - They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2699.6667
$ws.Range("I71").Value = 2700
$ws.Range("J71").Value = 2699.5
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 13497.5
$ws.Range("M71").Value = -9756
$ws.Range("N71").Value = -20985.5

# row 74 - Overall, We Blend In
$ws.Range("H74").Value = 31282.834
$ws.Range("I74").Value = 3848.5
$ws.Range("J74").Value = 45000
$ws.Range("K74").Value = 3848.5
$ws.Range("L74").Value = 45000
$ws.Range("M74").Value = -2850.5
$ws.Range("N74").Value = -46996

# row 77 - Eviction Notice (L)
$ws.Range("H77").Value = 31282.834
$ws.Range("I77").Value = 3848.5
$ws.Range("J77").Value = 45000
$ws.Range("K77").Value = 11545.5
$ws.Range("L77").Value = 135000
$ws.Range("M77").Value = -6553.5
$ws.Range("N77").Value = -144984

# row 82 - Trainin' the Neck
$ws.Range("H82").Value = 2750
$ws.Range("I82").Value = 2750
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2750
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2389
$ws.Range("N82").ClearContents()

# row 85 - Training Is Only Skintight (L)
$ws.Range("H85").Value = 2750
$ws.Range("I85").Value = 2750
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2750
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1502
$ws.Range("N85").ClearContents()

# row 136 - Respect for Br'aax
$ws.Range("H136").Value = 1709.45
$ws.Range("I136").Value = 1638.125
$ws.Range("K136").Value = 4914.375
$ws.Range("M136").Value = -2364.375

$ws = $wb.Worksheets.Item("WVR")
# row 62 - Pride Up in Smoke
$ws.Range("H62").Value = 4291.1665
$ws.Range("J62").Value = 4749.25
$ws.Range("L62").Value = 4749.25
$ws.Range("N62").Value = -5997.25

# row 65 - Desperate for Diversionaries (L)
$ws.Range("H65").Value = 4291.1665
$ws.Range("J65").Value = 4749.25
$ws.Range("L65").Value = 23746.25
$ws.Range("N65").Value = -29986.25

# row 68 - What Not to Wear
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# row 71 - Appeal of Foreign Apparel (L)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# row 132 - Comfy Cabins
$ws.Range("H132").Value = 991.5
$ws.Range("I132").Value = 647
$ws.Range("J132").Value = 2178.111
$ws.Range("K132").Value = 1941
$ws.Range("L132").Value = 6534.333
$ws.Range("M132").Value = 589
$ws.Range("N132").Value = -11594.333
